$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 180, shifting existing rows 180:291 down to 181:292
$ws.Rows("180:180").Insert()

# Populate the newly inserted row 180 with the new weekly data entry
$ws.Cells.Item(180, 1).Value = 10
$ws.Cells.Item(180, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(180, 3).Value = "La Araucanía"
$ws.Cells.Item(180, 4).Value = 45126
$ws.Cells.Item(180, 5).Value = 9
$ws.Cells.Item(180, 6).Value = 100112012
$ws.Cells.Item(180, 7).Value = "Espinaca"
$ws.Cells.Item(180, 8).Value = "Sin especificar"
$ws.Cells.Item(180, 9).Value = "Primera"
$ws.Cells.Item(180, 10).Value = 30
$ws.Cells.Item(180, 11).Value = 8000
$ws.Cells.Item(180, 12).Value = 8000
$ws.Cells.Item(180, 13).Value = 8000
$ws.Cells.Item(180, 14).Value = "`$/docena de paquetes"
$ws.Cells.Item(180, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(180, 16).Value = 667
$ws.Cells.Item(180, 17).Value = 12
$ws.Cells.Item(180, 18).Value = "Hortaliza"
